$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 164633
$ws.Range("C4").Value = 155584
$ws.Range("C5").Value = 9049
$ws.Range("C8").Value = 64.98999999999999
